$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("D23").Value = "3G Cells, 4G"
$ws.Range("E23").Value = ": HW Alarm: nan"
$ws.Range("G23").Value = "Partial"
$ws.Range("K23").Value = "300:37"

# Row 39
$ws.Range("D39").Value = "3G Cells, 4G"
$ws.Range("E39").Value = ": HW Alarm: nan"
$ws.Range("G39").Value = "Partial"
$ws.Range("K39").Value = "271:58"

# Row 75
$ws.Range("D75").Value = "3G Cells, 4G"
$ws.Range("E75").Value = ": HW Alarm: nan"
$ws.Range("G75").Value = "Partial"
$ws.Range("K75").Value = "304:53"

# Row 84
$ws.Range("D84").Value = "3G, 4G"
$ws.Range("E84").Value = ": HW Alarm: nan"
$ws.Range("G84").Value = "Partial"
$ws.Range("K84").Value = "416:02"

# Row 111
$ws.Range("K111").Value = "34:48"

# Row 115
$ws.Range("K115").Value = "07:47"

# Row 123
$ws.Range("K123").Value = "196:01"

# Row 151
$ws.Range("D151").Value = "3G Cells, 4G"
$ws.Range("E151").Value = ": HW Alarm: nan<br>: HW Alarm: cell disabled<br>LTE: LTE Cells: 1<br>LTE: HW Alarm: VSWR<br>LTE: VSWR (ANTE)<br>WCDMA: WCDMA Cells: 2<br>WCDMA: HW Alarm: VSWR<br>WCDMA: VSWR (ANTE)<br>WCDMA: WCDMA Cells: 1"
$ws.Range("G151").Value = "Partial"
$ws.Range("K151").Value = "639:27"
